$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in / clear a few individual cells first (row numbers unaffected by later row deletes)
$ws.Range("D19").Value = -15.5
$ws.Range("D21").Value = $null
$ws.Range("D23").Value = -13.9

# Delete entire rows for "RM 232" (row 26) and "SC 92" (row 28).
# Delete the lower row first so the earlier row index stays valid.
$ws.Rows(28).Delete()
$ws.Rows(26).Delete()

# After the two row deletions, the data that used to be on rows 27/29/31/35 is now on rows 26/27/29/33
$ws.Range("B26").Value = $null
$ws.Range("B27").Value = -20.4
$ws.Range("D27").Value = $null
$ws.Range("B29").Value = $null
$ws.Range("D33").Value = -14.1
